# feat: add 2022-Q1 data
#
# The existing "总计" (Total) sheet becomes the new "2022-Q1" sheet (keeps
# its sheetId/rId), and a brand-new "总计" sheet is appended after it that
# aggregates the historical table plus the new 2022-Q1 row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: turn the current "总计" sheet into the new "2022-Q1" holdings
# sheet (same position-relative sheetId assignment Excel would use when
# you rename the existing tab and then add a fresh one after it).
# ---------------------------------------------------------------------
$q1 = $wb.Worksheets.Item("总计")
$q1.Cells.Clear()
$q1.Name = "2022-Q1"

# headers (row 1, columns B:H) - bold/centered header style lives on the
# "2021-Q4" sheet's header row; we copy formats only, after values.
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

function Set-TextValue($range, [string]$text) {
    # Force a numeric-looking string to be stored as text (keeps leading
    # zeros / decimal formatting identical to the source value), then
    # drop the "Text" number-format style iron_native applies so the
    # cell ends up with the workbook's default (no) style, same as the
    # source data.
    $range.Value = "'" + $text
    $range.ClearFormats()
}

# row 2
$q1.Range("A2").Value = 0
Set-TextValue $q1.Range("B2") "377010"
$q1.Range("C2").Value = "上投摩根阿尔法混合"
Set-TextValue $q1.Range("D2") "13.21"
Set-TextValue $q1.Range("E2") "89.24"
Set-TextValue $q1.Range("F2") "7.11"
Set-TextValue $q1.Range("G2") "0.9392"
$q1.Range("H2").Value = 1

# row 3
$q1.Range("A3").Value = 1
Set-TextValue $q1.Range("B3") "000457"
$q1.Range("C3").Value = "上投摩根核心成长"
Set-TextValue $q1.Range("D3") "12.84"
Set-TextValue $q1.Range("E3") "87.63"
Set-TextValue $q1.Range("F3") "6.67"
Set-TextValue $q1.Range("G3") "0.8564"
$q1.Range("H3").Value = 2

# row 4
$q1.Range("A4").Value = 2
Set-TextValue $q1.Range("B4") "000892"
$q1.Range("C4").Value = "九泰天宝灵活配置混合A"
Set-TextValue $q1.Range("D4") "0.07"
Set-TextValue $q1.Range("E4") "90.81"
Set-TextValue $q1.Range("F4") "4.95"
Set-TextValue $q1.Range("G4") "0.0035"
$q1.Range("H4").Value = 4

# row 5
$q1.Range("A5").Value = 3
Set-TextValue $q1.Range("B5") "002028"
$q1.Range("C5").Value = "九泰天宝灵活配置混合C"
Set-TextValue $q1.Range("D5") "0.00"
Set-TextValue $q1.Range("E5") "90.81"
Set-TextValue $q1.Range("F5") "4.95"
$q1.Range("G5").Value = 0
$q1.Range("H5").Value = 4

# Re-apply the shared header/index style (bordered, bold, centered) onto
# the header row and the A-column index cells, matching every other
# quarter sheet in the workbook.
$styleSrc = $wb.Worksheets.Item("2021-Q4")
$styleSrc.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)
$styleSrc.Range("A2:A5").Copy()
$q1.Range("A2:A5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# Step 2: append a fresh "总计" sheet after "2022-Q1" with the updated
# roll-up table (new 2022-Q1 row on top, everything else shifted down).
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $q1)
$total.Name = "总计"

$total.Range("B1").Value = "日期"
$total.Range("C1").Value = "持有数量(只)"
$total.Range("D1").Value = "持有市值(亿元)"

$totalRows = @(
    @(0, "2022-Q1", 4, 1.8),
    @(1, "2021-Q4", 5, 2.59),
    @(2, "2021-Q3", 6, 3.4),
    @(3, "2021-Q2", 8, 2.93),
    @(4, "2021-Q1", 8, 2.86),
    @(5, "2020-Q4", 9, 3.27)
)

$r = 2
foreach ($row in $totalRows) {
    $total.Range("A$r").Value = $row[0]
    $total.Range("B$r").Value = $row[1]
    $total.Range("C$r").Value = $row[2]
    $total.Range("D$r").Value = $row[3]
    $r = $r + 1
}

$styleSrc.Range("B1:D1").Copy()
$total.Range("B1:D1").PasteSpecial(-4122)
# "2021-Q3" is the sheet with 6 rows of index cells (A2:A7), matching the
# roll-up table's new row count exactly.
$indexStyleSrc = $wb.Worksheets.Item("2021-Q3")
$indexStyleSrc.Range("A2:A7").Copy()
$total.Range("A2:A7").PasteSpecial(-4122)
$excel.CutCopyMode = $false
